$d = $word.ActiveDocument

# --- 1) Locate the "Library(" call (the one right before "BHH2)") and split
#        it into two runs "l" + "ibrary(", also correcting the stray leading
#        capital L -> l (R's library() function is lower-case).
#
#        A plain Range.Text assignment just rewrites the text in place without
#        creating a new run boundary, so to force the run split we briefly
#        anchor the (uniquely named, singleton) _GoBack bookmark at the split
#        point. The run split persists even after the bookmark is relocated
#        in step 2 below.
$rLib = $d.Content
$foundLib = $rLib.Find.Execute("Library(BHH2)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundLib) {
  throw "Could not find 'Library(BHH2)' in the document"
}
$libStart = $rLib.Start

$rL = $d.Range($libStart, $libStart + 1)
$rL.Text = "l"

$rLibSplit = $d.Range($libStart + 1, $libStart + 1)
$d.Bookmarks.Add("_GoBack", $rLibSplit) | Out-Null

# --- 2) Locate "these" inside "...produce these plots in R..." and split it
#        into "the" + "se" (no text change, just a run boundary). This is
#        where the user's cursor/last-edit landed, so the _GoBack bookmark
#        (Word's "last edit position" marker, previously sitting elsewhere in
#        the document) ends up here -- its final resting place.
$rThese = $d.Content
$foundThese = $rThese.Find.Execute("produce these plots", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundThese) {
  throw "Could not find 'produce these plots' in the document"
}
$theseStart = $rThese.Start + "produce ".Length

$rSplit = $d.Range($theseStart + 3, $theseStart + 3)
$d.Bookmarks.Add("_GoBack", $rSplit) | Out-Null
